# Update the "Trading History" sheet with 3 new Buy trades.
#
# The existing row 5 (date 2025-10-29 / serial 45959, qty 50 @ 138.6,
# CN#252607497001) moves down to row 8 unchanged, and three new trade
# rows are written above it at rows 5-7 (newest trade first):
#   row 5 -> 2025-11-10 (45971), qty 20  @ 144.53, CN#252607962278
#   row 6 -> 2025-10-31 (45961), qty 30  @ 143.12, CN#252607606298
#   row 7 -> 2025-10-30 (45960), qty 50  @ 139.52, CN#252607551490
#   row 8 -> 2025-10-29 (45959), qty 50  @ 138.6,  CN#252607497001 (orig row 5)
#
# Every row keeps the same column layout as the existing data rows:
#   A=DATE  B=EXCH  C=ACTION  D=QTY  E=PRICE  F=COST  G=REMARKS
#   I=ADD CHRG   J=Current Price (formula =Index!$C$2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

$dateFormat = "yyyy-mm-dd h:mm:ss"

$rows = @(
    @{ r = 5; date = 45971; exch = "NSE"; action = "Buy"; qty = 20; price = 144.53; cost = 2910.6; remarks = "CN#252607962278" },
    @{ r = 6; date = 45961; exch = "NSE"; action = "Buy"; qty = 30; price = 143.12; cost = 4313.6; remarks = "CN#252607606298" },
    @{ r = 7; date = 45960; exch = "NSE"; action = "Buy"; qty = 50; price = 139.52; cost = 6996;   remarks = "CN#252607551490" },
    @{ r = 8; date = 45959; exch = "NSE"; action = "Buy"; qty = 50; price = 138.6;  cost = 6950;   remarks = "CN#252607497001" }
)

foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value = $row.date
    $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 2).Value = $row.exch
    $ws.Cells.Item($r, 3).Value = $row.action
    $ws.Cells.Item($r, 4).Value = $row.qty
    $ws.Cells.Item($r, 5).Value = $row.price
    $ws.Cells.Item($r, 6).Value = $row.cost
    $ws.Cells.Item($r, 7).Value = $row.remarks
    $ws.Cells.Item($r, 9).Value = 20
    $ws.Cells.Item($r, 10).Formula = "=Index!`$C`$2"
}
